# Update gh-pages output (generated data refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Sheet1): update "想去人数" (F column) counts ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value2 = 13233
$ws1.Cells.Item(3, 6).Value2 = 319
$ws1.Cells.Item(5, 6).Value2 = 213
$ws1.Cells.Item(6, 6).Value2 = 425
$ws1.Cells.Item(7, 6).Value2 = 1285
$ws1.Cells.Item(8, 6).Value2 = 121

# --- Sheet "演出" (Sheet2): drop the oldest event, shift remaining rows up, drop last row ---
$ws2 = $wb.Worksheets.Item("演出")
# Force column B to Text first so the "yyyy-mm-dd" strings aren't auto-converted
# into date serials when copied; restore the default ("Normal") style afterwards
# so no stray number-format is left behind on the cells.
$ws2.Range("B2:B4").NumberFormat = "@"
$ws2.Range("B2:I4").Value2 = $ws2.Range("B3:I5").Value2
$ws2.Range("B2:B4").Style = "Normal"
$ws2.Rows.Item(5).Delete()

# --- Sheet "全部类型" (Sheet4): drop the oldest event, shift remaining rows up, drop last row ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("B2:B11").NumberFormat = "@"
$ws4.Range("B2:I11").Value2 = $ws4.Range("B3:I12").Value2
$ws4.Range("B2:B11").Style = "Normal"
$ws4.Rows.Item(12).Delete()
